$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 8900
$ws.Range("I111").Value = 10000
$ws.Range("J111").Value = 7800
$ws.Range("K111").Value = 30000
$ws.Range("L111").Value = 23400
$ws.Range("M111").Value = -26933
$ws.Range("N111").Value = -29534

$ws.Range("H137").Value = 1545.7954
$ws.Range("I137").Value = 1430.3077
$ws.Range("J137").Value = 1712.6111
$ws.Range("K137").Value = 4290.9231
$ws.Range("L137").Value = 5137.8333
$ws.Range("M137").Value = -1740.9231
$ws.Range("N137").Value = -10237.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1230.6471
$ws.Range("I2").Value = 1272.2
$ws.Range("J2").Value = 1171.2858
$ws.Range("K2").Value = 1272.2
$ws.Range("L2").Value = 1171.2858
$ws.Range("M2").Value = -1159.2
$ws.Range("N2").Value = -1397.2858

$ws.Range("H32").Value = 4287.62
$ws.Range("I32").Value = 3708.5
$ws.Range("J32").Value = 9499.7
$ws.Range("K32").Value = 3708.5
$ws.Range("L32").Value = 9499.7
$ws.Range("M32").Value = -3421.5
$ws.Range("N32").Value = -10073.7

$ws.Range("H102").Value = 9804973
$ws.Range("I102").Value = 10417734
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 10417734
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -10416112
$ws.Range("N102").Value = -4044

$ws.Range("H116").Value = 1230.6471
$ws.Range("I116").Value = 1272.2
$ws.Range("J116").Value = 1171.2858
$ws.Range("K116").Value = 1272.2
$ws.Range("L116").Value = 1171.2858
$ws.Range("M116").Value = 1021.8
$ws.Range("N116").Value = -5759.2858

$ws.Range("H132").Value = 2312.468
$ws.Range("I132").Value = 1949.8889
$ws.Range("J132").Value = 2801.95
$ws.Range("K132").Value = 5849.6667
$ws.Range("L132").Value = 8405.849999999999
$ws.Range("M132").Value = -3319.6667
$ws.Range("N132").Value = -13465.85

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1230.6471
$ws.Range("I3").Value = 1272.2
$ws.Range("J3").Value = 1171.2858
$ws.Range("K3").Value = 1272.2
$ws.Range("L3").Value = 1171.2858
$ws.Range("M3").Value = -1158.2
$ws.Range("N3").Value = -1399.2858

$ws.Range("H20").Value = 2985.3333
$ws.Range("J20").Value = 3296.6667
$ws.Range("L20").Value = 3296.6667
$ws.Range("N20").Value = -3790.6667

$ws.Range("H132").Value = 32140
$ws.Range("J132").Value = 32140
$ws.Range("L132").Value = 32140
$ws.Range("N132").Value = -42260

$ws.Range("H134").Value = 991.5217
$ws.Range("I134").Value = 847.9048
$ws.Range("K134").Value = 2543.7144
$ws.Range("M134").Value = -8.714399999999841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1532.8684
$ws.Range("I31").Value = 1541.8649
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 1541.8649
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = -1246.8649
$ws.Range("N31").Value = -1790

$ws.Range("H34").Value = 1532.8684
$ws.Range("I34").Value = 1541.8649
$ws.Range("J34").Value = 1200
$ws.Range("K34").Value = 1541.8649
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -1339.8649
$ws.Range("N34").Value = -1604

$ws.Range("H86").Value = 3367893
$ws.Range("I86").Value = 4466723
$ws.Range("K86").Value = 4466723
$ws.Range("M86").Value = -4465600

$ws.Range("H89").Value = 3367893
$ws.Range("I89").Value = 4466723
$ws.Range("K89").Value = 22333615
$ws.Range("M89").Value = -22327999

$ws.Range("H114").Value = 31995
$ws.Range("J114").Value = 31995
$ws.Range("L114").Value = 31995
$ws.Range("N114").Value = -40673

$ws.Range("H132").Value = 1735.1082
$ws.Range("I132").Value = 1412.1852
$ws.Range("K132").Value = 4236.5556
$ws.Range("M132").Value = -1706.5556

$ws.Range("H134").Value = 12821631
$ws.Range("I134").Value = 1034.2858
$ws.Range("K134").Value = 3102.8574
$ws.Range("M134").Value = -567.8574000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 117.63158
$ws.Range("I12").Value = 212.8
$ws.Range("J12").Value = 83.64286
$ws.Range("K12").Value = 638.4000000000001
$ws.Range("L12").Value = 250.92858
$ws.Range("M12").Value = -465.4000000000001
$ws.Range("N12").Value = -596.92858

$ws.Range("H125").Value = 2765
$ws.Range("I125").Value = 1353.3334
$ws.Range("J125").Value = 7000
$ws.Range("K125").Value = 4060.0002
$ws.Range("L125").Value = 21000
$ws.Range("M125").Value = 859.9998
$ws.Range("N125").Value = -30840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40912350
$ws.Range("I70").Value = 35717764
$ws.Range("J70").Value = 50002876
$ws.Range("K70").Value = 35717764
$ws.Range("L70").Value = 50002876
$ws.Range("M70").Value = -35717494
$ws.Range("N70").Value = -50003416

$ws.Range("H73").Value = 40912350
$ws.Range("I73").Value = 35717764
$ws.Range("J73").Value = 50002876
$ws.Range("K73").Value = 35717764
$ws.Range("L73").Value = 50002876
$ws.Range("M73").Value = -35716828
$ws.Range("N73").Value = -50004748

$ws.Range("H126").Value = 1896.7646
$ws.Range("I126").Value = 1539
$ws.Range("J126").Value = 2299.25
$ws.Range("K126").Value = 4617
$ws.Range("L126").Value = 6897.75
$ws.Range("M126").Value = -2147
$ws.Range("N126").Value = -11837.75

$ws.Range("H132").Value = 3377.1738
$ws.Range("I132").Value = 3122.6875
$ws.Range("J132").Value = 3958.8572
$ws.Range("K132").Value = 9368.0625
$ws.Range("L132").Value = 11876.5716
$ws.Range("M132").Value = -6838.0625
$ws.Range("N132").Value = -16936.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2373
$ws.Range("I7").Value = 1594.8
$ws.Range("K7").Value = 1594.8
$ws.Range("M7").Value = -1482.8

$ws.Range("H40").Value = 5494.4287
$ws.Range("I40").Value = 3071.625
$ws.Range("J40").Value = 8724.833
$ws.Range("K40").Value = 3071.625
$ws.Range("L40").Value = 8724.833
$ws.Range("M40").Value = -2935.625
$ws.Range("N40").Value = -8996.833

$ws.Range("H100").Value = 1082.1538
$ws.Range("I100").Value = 1024
$ws.Range("J100").Value = 1150
$ws.Range("K100").Value = 1024
$ws.Range("L100").Value = 1150
$ws.Range("M100").Value = -483
$ws.Range("N100").Value = -2232

$ws.Range("H126").Value = 2373
$ws.Range("I126").Value = 1594.8
$ws.Range("K126").Value = 4784.4
$ws.Range("M126").Value = -2314.4

$ws.Range("H132").Value = 3682.7144
$ws.Range("I132").Value = 4254.3335
$ws.Range("J132").Value = 3254
$ws.Range("K132").Value = 12763.0005
$ws.Range("L132").Value = 9762
$ws.Range("M132").Value = -10233.0005
$ws.Range("N132").Value = -14822

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 29166.334
$ws.Range("J98").Value = 29166.334
$ws.Range("L98").Value = 29166.334
$ws.Range("N98").Value = -35156.334

$ws.Range("H100").Value = 355.54544
$ws.Range("I100").Value = 354.57144
$ws.Range("J100").Value = 357.25
$ws.Range("K100").Value = 709.14288
$ws.Range("L100").Value = 714.5
$ws.Range("M100").Value = -168.14288
$ws.Range("N100").Value = -1796.5

$ws.Range("H115").Value = 36708.363
$ws.Range("I115").Value = 25000
$ws.Range("J115").Value = 37879.2
$ws.Range("K115").Value = 25000
$ws.Range("L115").Value = 37879.2
$ws.Range("M115").Value = -23433
$ws.Range("N115").Value = -41013.2

$ws.Range("H126").Value = 66667096
$ws.Range("I126").Value = 66667096
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 200001288
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -199998818
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2291
$ws.Range("I132").Value = 2127.7
$ws.Range("J132").Value = 2757.5715
$ws.Range("K132").Value = 6383.099999999999
$ws.Range("L132").Value = 8272.7145
$ws.Range("M132").Value = -3853.099999999999
$ws.Range("N132").Value = -13332.7145

$ws.Range("H136").Value = 1410.88
$ws.Range("I136").Value = 1250.421
$ws.Range("J136").Value = 1919
$ws.Range("K136").Value = 3751.263
$ws.Range("L136").Value = 5757
$ws.Range("M136").Value = -1201.263
$ws.Range("N136").Value = -10857
